$d = $word.ActiveDocument

# --- Change 1: merge the "3" + ".  " runs into a single run "3.  " -------
# (leaves the following "Configuración Técnica (MkDocs):" run, with its
# underline formatting, untouched, since the match stops right at the
# run boundary). MatchWildcards is off so the "." is matched literally.
$d.Content.Find.Execute("3.  ", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "3.  ", 2, $false) | Out-Null

# --- Change 2: remove the stray blank paragraphs + the "Paso 1..Conclusiones"
# draft/outline block that followed the "GitHub Pages" paragraph (the one
# carrying the _GoBack bookmark). After the deletion that paragraph is
# directly followed by the two pre-existing trailing empty paragraphs. ----
$goBackPara = $d.Bookmarks("_GoBack").Range.Paragraphs(1)
$startPara = $goBackPara.Next()

$conclusionesPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq ("Conclusiones" + [char]13)) {
        $conclusionesPara = $p
    }
}

$delRange = $d.Range($startPara.Range.Start, $conclusionesPara.Range.End)
$delRange.Delete()
